# Update the "F" column (想去人数 / want-to-go count) values per the source diff.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 6693
$ws.Range("F3").Value = 796
$ws.Range("F4").Value = 1100
$ws.Range("F5").Value = 130
$ws.Range("F6").Value = 702
$ws.Range("F7").Value = 220
$ws.Range("F8").Value = 12
$ws.Range("F9").Value = 1081
$ws.Range("F10").Value = 842
$ws.Range("F11").Value = 1007
$ws.Range("F12").Value = 1320
$ws.Range("F13").Value = 37
$ws.Range("F14").Value = 110
$ws.Range("F15").Value = 527
$ws.Range("F18").Value = 363
$ws.Range("F19").Value = 1054
$ws.Range("F20").Value = 1477
$ws.Range("F22").Value = 205
$ws.Range("F23").Value = 447
$ws.Range("F24").Value = 437
$ws.Range("F26").Value = 4
$ws.Range("F27").Value = 1122
$ws.Range("F28").Value = 246
$ws.Range("F29").Value = 2341
$ws.Range("F31").Value = 1199
$ws.Range("F32").Value = 433
$ws.Range("F34").Value = 3814

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 1033
$ws.Range("F19").Value = 4119
$ws.Range("F22").Value = 22
$ws.Range("F27").Value = 105
$ws.Range("F29").Value = 223

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 1243
$ws.Range("F5").Value = 1624
$ws.Range("F7").Value = 135
$ws.Range("F8").Value = 939

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1243
$ws.Range("F3").Value = 1624
$ws.Range("F5").Value = 135
$ws.Range("F6").Value = 939
$ws.Range("F9").Value = 6693
$ws.Range("F12").Value = 796
$ws.Range("F14").Value = 130
$ws.Range("F15").Value = 702
$ws.Range("F16").Value = 1081
$ws.Range("F17").Value = 842
$ws.Range("F20").Value = 1007
$ws.Range("F22").Value = 1320
$ws.Range("F23").Value = 37
$ws.Range("F24").Value = 527
$ws.Range("F27").Value = 363
$ws.Range("F28").Value = 1477
$ws.Range("F30").Value = 447
$ws.Range("F31").Value = 437
$ws.Range("F34").Value = 1122
$ws.Range("F35").Value = 246
$ws.Range("F36").Value = 105
$ws.Range("F38").Value = 2341
$ws.Range("F39").Value = 223
$ws.Range("F45").Value = 1199
$ws.Range("F46").Value = 433
$ws.Range("F48").Value = 3814
